$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the data values currently in row 2 before we overwrite row 1.
$bankName = $ws.Range("B2").Value()
$depositType = $ws.Range("C2").Value()
$currency = $ws.Range("D2").Value()
$owner = $ws.Range("E2").Value()
$total = $ws.Range("F2").Value()

# Move the data values up into row 1, replacing the header text.
$ws.Range("B1").Value = $bankName
$ws.Range("C1").Value = $depositType
$ws.Range("D1").Value = $currency
$ws.Range("E1").Value = $owner
$ws.Range("F1").Value = $total

# Clear column A (the "12" land-portion style value is no longer used) and
# delete the old row 2 entirely, so only row 1 (B1:F1) remains populated.
$ws.Range("A1:A2").ClearContents()
$ws.Rows.Item(2).Delete()
